$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Go get analog single ended 2:1 multiplexer SPDT low level takes 5V"
#    paragraph -> make bold (both runs + paragraph mark run props)
# ------------------------------------------------------------------
$pBold = $d.Paragraphs.Item(13)
if ($pBold.Range.Text -notmatch "Go get analog single ended") {
    throw "paragraph 13 text mismatch: $($pBold.Range.Text)"
}
$pBold.Range.Font.Bold = $true

# ------------------------------------------------------------------
# 2) Remove the _GoBack bookmark currently at the end of the
#    "want power off switch..." paragraph.
# ------------------------------------------------------------------
$pSwitch = $d.Paragraphs.Item(11)
if ($pSwitch.Range.Text -notmatch "want power off switch") {
    throw "paragraph 11 text mismatch: $($pSwitch.Range.Text)"
}
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# ------------------------------------------------------------------
# 3) "The input trigger voltage… is trigger input voltage threshold"
#    paragraph -> split into 5 runs by inserting new text, and move
#    the _GoBack bookmark to the (true) end of this paragraph.
# ------------------------------------------------------------------
$pTrigger = $d.Paragraphs.Item(19)
if ($pTrigger.Range.Text -notmatch "The input trigger voltage") {
    throw "paragraph 19 text mismatch: $($pTrigger.Range.Text)"
}
$pStart = $pTrigger.Range.Start

# "The input trigger voltage" | " on the data sheet of the timer" | "… is" | " actually " | " trigger input voltage threshold"
$newPiece1 = " on the data sheet of the timer"
$newPiece2 = " actually "

# insert the new text pieces first (this merges into the existing runs;
# we split them into their own <w:r> elements afterwards)
$insertPos1 = $pStart + 25   # right after "The input trigger voltage"
$r1 = $d.Range($insertPos1, $insertPos1)
$r1.InsertAfter($newPiece1)

$insertPos2 = $pStart + 25 + $newPiece1.Length + 4   # right after "...timer… is"
$r2 = $d.Range($insertPos2, $insertPos2)
$r2.InsertAfter($newPiece2)

# split the paragraph into 5 pieces (at descending offsets so earlier
# offsets remain valid), then re-join by deleting the paragraph marks;
# this preserves each piece as its own <w:r> instead of Word silently
# re-merging adjacent same-formatted text back into one run.
$boundaries = @(70, 60, 56, 25)
foreach ($b in $boundaries) {
    $pos = $pStart + $b
    $rSplit = $d.Range($pos, $pos)
    $rSplit.InsertParagraphAfter()
}
for ($k = 0; $k -lt 4; $k++) {
    $pMerge = $d.Paragraphs.Item(19)
    $endRange = $d.Range($pMerge.Range.End - 1, $pMerge.Range.End)
    $endRange.Delete()
}

$pTrigger = $d.Paragraphs.Item(19)
$expectedText = "The input trigger voltage" + $newPiece1 + [char]0x2026 + " is" + $newPiece2 + " trigger input voltage threshold"
if ($pTrigger.Range.Text -ne $expectedText) {
    throw "unexpected paragraph 19 text after split: [$($pTrigger.Range.Text)]"
}

# Move the _GoBack bookmark to the true end of this paragraph (right
# before its paragraph mark). A zero-width bookmark placed exactly at
# a paragraph boundary is mis-anchored, so temporarily append a filler
# character, anchor the bookmark before it, then remove the filler.
$endOfText = $pTrigger.Range.Duplicate
$endOfText.Collapse(0)
$endOfText.MoveEnd(1, -1)
$fillerPos = $endOfText.Start
$endOfText.InsertAfter("X")

$bmRange = $d.Range($fillerPos, $fillerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$fillerRange = $d.Range($fillerPos, $fillerPos + 1)
$fillerRange.Delete()

# ------------------------------------------------------------------
# 4) "ATV transmitter will be off the 5V supply " paragraph -> bold + red
# ------------------------------------------------------------------
$pATV = $d.Paragraphs.Item(20)
if ($pATV.Range.Text -notmatch "ATV transmitter will be off the 5V supply") {
    throw "paragraph 20 text mismatch: $($pATV.Range.Text)"
}
$pATV.Range.Font.Bold = $true
$pATV.Range.Font.Color = 255
